$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 03:04:24"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 03:04:19"
$wsZhCn.Range("K2").Value = "2016-08-28 03:04:38"

# de-de sheet: "Correspond Handoff Datetime" (shares the same text as Overview!G2)
# and "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 03:04:24"
$wsDeDe.Range("K2").Value = "2016-08-28 03:04:44"
